# Adds the 2024/12/16 column (CU, index 99) to the "合成確率" sheet,
# matching the style convention already used by every other date column:
#   s=1 -> no fill      (value >= 140)
#   s=2 -> yellow fill  (value <  125)
#   s=3 -> light-blue fill (125 <= value < 140)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("合成確率")

$newCol = 99   # column CU

# New column gets the same width (12 chars) as every other data column.
# ColumnWidth is in "characters"; this value round-trips to xml width="12"
# for this workbook's default font, same as columns 1..98.
$ws.Columns.Item($newCol).ColumnWidth = 11.1666666666667

# --- Header cell CU1: "2024/12/16" as literal text, not an auto-parsed date ---
$cu1 = $ws.Cells.Item(1, $newCol)
$cu1.NumberFormat = "@"
$cu1.Value = "2024/12/16"
# Re-apply the same (fill-less) style used by every other header cell, without
# disturbing the "stored as text" flag that NumberFormat "@" gave us above.
$ws.Cells.Item(1, 1).Copy()
$cu1.PasteSpecial(-4122)

# --- Reference cells already on the sheet carrying each of the three styles ---
$styleNormal = $ws.Cells.Item(2, 1)    # s=1 no fill
$styleLow    = $ws.Cells.Item(2, 4)    # s=2 yellow fill
$styleMid    = $ws.Cells.Item(2, 14)   # s=3 light-blue fill

# --- New data values for 2024/12/16, one per machine (row 2..53) ---
$data = @(
    @{Row=2; Value=130.8; Style=3},
    @{Row=3; Value=192.3; Style=1},
    @{Row=4; Value=140; Style=1},
    @{Row=5; Value=211; Style=1},
    @{Row=6; Value=173.1; Style=1},
    @{Row=7; Value=144.9; Style=1},
    @{Row=8; Value=127.8; Style=3},
    @{Row=9; Value=183.5; Style=1},
    @{Row=10; Value=160.4; Style=1},
    @{Row=11; Value=154; Style=1},
    @{Row=12; Value=129.4; Style=3},
    @{Row=13; Value=136.2; Style=3},
    @{Row=14; Value=465.5; Style=1},
    @{Row=15; Value=146.7; Style=1},
    @{Row=16; Value=117.7; Style=2},
    @{Row=17; Value=173.3; Style=1},
    @{Row=18; Value=170.8; Style=1},
    @{Row=19; Value=170.8; Style=1},
    @{Row=20; Value=148.1; Style=1},
    @{Row=21; Value=138.7; Style=3},
    @{Row=22; Value=156; Style=1},
    @{Row=23; Value=393.5; Style=1},
    @{Row=24; Value=168.2; Style=1},
    @{Row=25; Value=133.5; Style=3},
    @{Row=26; Value=204; Style=1},
    @{Row=27; Value=152.2; Style=1},
    @{Row=28; Value=131.8; Style=3},
    @{Row=29; Value=173.4; Style=1},
    @{Row=30; Value=177.8; Style=1},
    @{Row=31; Value=163.6; Style=1},
    @{Row=32; Value=195; Style=1},
    @{Row=33; Value=161.7; Style=1},
    @{Row=34; Value=136.5; Style=3},
    @{Row=35; Value=164.7; Style=1},
    @{Row=36; Value=114.6; Style=2},
    @{Row=37; Value=206.8; Style=1},
    @{Row=38; Value=132.7; Style=3},
    @{Row=39; Value=153.9; Style=1},
    @{Row=40; Value=180.9; Style=1},
    @{Row=41; Value=188.9; Style=1},
    @{Row=42; Value=133.9; Style=3},
    @{Row=43; Value=165.5; Style=1},
    @{Row=44; Value=206.4; Style=1},
    @{Row=45; Value=141.8; Style=1},
    @{Row=46; Value=143.1; Style=1},
    @{Row=47; Value=170.7; Style=1},
    @{Row=48; Value=133.6; Style=3},
    @{Row=49; Value=163.6; Style=1},
    @{Row=50; Value=417; Style=1},
    @{Row=51; Value=146; Style=1},
    @{Row=52; Value=238.1; Style=1},
    @{Row=53; Value=146.2; Style=1}
)

foreach ($item in $data) {
    $cell = $ws.Cells.Item($item.Row, $newCol)
    $cell.Value = $item.Value

    if ($item.Style -eq 2) {
        $styleLow.Copy()
    } elseif ($item.Style -eq 3) {
        $styleMid.Copy()
    } else {
        $styleNormal.Copy()
    }
    $cell.PasteSpecial(-4122)
}
